# Update inventory figures on Sheet1 and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sold (B) / Price (C) updates
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 28

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 20

$ws.Range("B4").Value = 35
$ws.Range("C4").Value = 49

$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 152

$ws.Range("B6").Value = 15
$ws.Range("C6").Value = 12

$ws.Range("B8").Value = 34
$ws.Range("C8").Value = 444

$ws.Range("B10").Value = 50
$ws.Range("C10").Value = 40

# Move the active selection to C10, matching the saved view state.
$ws.Range("C10").Select()
